$wb = $excel.ActiveWorkbook

# --- BME_BCCW sheet: update the "partial overhead" formulas -------------
# Column B (actual_partial_oh) used to be total_exp - labour - contracts - parts.
# Contracts/parts are now dropped so duplicate asset descriptions no longer
# throw the calc off: it's simply total_exp - labour_exp.
# Column G (budgeted_partial_oh) mirrors the same change for the budget side.
$ws2 = $wb.Worksheets.Item("BME_BCCW")
$ws2.Range("B2").Formula = "=C2-D2"
$ws2.Range("B3:B6").Formula = "=C3-D3"
$ws2.Range("G2").Formula = "=H2-I2"
$ws2.Range("G3:G6").Formula = "=H3-I3"

# --- BME_BCC sheet: same formula fix -------------------------------------
$ws3 = $wb.Worksheets.Item("BME_BCC")
$ws3.Range("B2").Formula = "=C2-D2"
$ws3.Range("B3:B6").Formula = "=C3-D3"
$ws3.Range("G2").Formula = "=H2-I2"
$ws3.Range("G3:G6").Formula = "=H3-I3"

# --- restore each sheet's own selection/scroll state ---------------------
$ws2.Range("D2").Select()
$ws3.Range("G28").Select()

# README ends up the active/visible sheet and cell when the file is saved.
$ws1 = $wb.Worksheets.Item("README")
$ws1.Activate()
$ws1.Range("E16").Select()
